$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values to reflect the new/expanded dataset ---
$ws.Range("T1").Value = 18
$ws.Range("L2").Value = '1Pair-B'
$ws.Range("M2").Value = '2Pairs-A'
$ws.Range("N2").Value = '2Pairs-B'
$ws.Range("O2").Value = '3Pairs-A'
$ws.Range("P2").Value = '3Pairs-B'
$ws.Range("Q2").Value = '3Pairs-C'
$ws.Range("R2").Value = '4Pairs'
$ws.Range("S2").Value = '5A4F'
$ws.Range("T2").Value = 'MaxUnique'
$ws.Range("C3").Value = 0.9507204610951009
$ws.Range("D3").Value = 1.045842939481268
$ws.Range("E3").Value = 1.001059077809798
$ws.Range("F3").Value = 0.9507204610951009
$ws.Range("G3").Value = 1.008364553314121
$ws.Range("H3").Value = 1.022853025936599
$ws.Range("I3").Value = 0.9885951008645533
$ws.Range("J3").Value = 1.045842939481268
$ws.Range("K3").Value = 0.9507204610951009
$ws.Range("L3").Value = 1.001059077809798
$ws.Range("M3").Value = 1.023451008645533
$ws.Range("N3").Value = 1.023451008645533
$ws.Range("O3").Value = 1.018422190201729
$ws.Range("P3").Value = 0.9992074927953891
$ws.Range("Q3").Value = 0.9992074927953891
$ws.Range("R3").Value = 0.987085734870317
$ws.Range("S3").Value = 0.987085734870317
$ws.Range("T3").Value = 1.00290585975024
$ws.Range("C4").Value = 1.001276867667034
$ws.Range("D4").Value = 0.9927242129093838
$ws.Range("E4").Value = 0.9929987323325614
$ws.Range("F4").Value = 1.001276867667034
$ws.Range("G4").Value = 0.995181361977499
$ws.Range("H4").Value = 0.9876187527966136
$ws.Range("I4").Value = 0.9949752787408748
$ws.Range("J4").Value = 0.9927242129093838
$ws.Range("K4").Value = 1.001276867667034
$ws.Range("L4").Value = 0.9929987323325614
$ws.Range("M4").Value = 0.9928614726209726
$ws.Range("N4").Value = 0.9928614726209726
$ws.Range("O4").Value = 0.9936347690731481
$ws.Range("P4").Value = 0.995666604302993
$ws.Range("Q4").Value = 0.995666604302993
$ws.Range("R4").Value = 0.9970691701440033
$ws.Range("S4").Value = 0.9970691701440033
$ws.Range("T4").Value = 0.9941292010706612
$ws.Range("C5").Value = 1.003714561987474
$ws.Range("D5").Value = 0.9948710600980648
$ws.Range("E5").Value = 0.9909624085808176
$ws.Range("F5").Value = 1.003714561987474
$ws.Range("G5").Value = 0.9970387717646837
$ws.Range("H5").Value = 0.9819478282985443
$ws.Range("I5").Value = 0.9942763554956564
$ws.Range("J5").Value = 0.9948710600980648
$ws.Range("K5").Value = 1.003714561987474
$ws.Range("L5").Value = 0.9909624085808176
$ws.Range("M5").Value = 0.9929167343394412
$ws.Range("N5").Value = 0.9929167343394412
$ws.Range("O5").Value = 0.994290746814522
$ws.Range("P5").Value = 0.9965160102221186
$ws.Range("Q5").Value = 0.9965160102221188
$ws.Range("R5").Value = 0.9983156481634574
$ws.Range("S5").Value = 0.9983156481634574
$ws.Range("T5").Value = 0.9938018310375401
$ws.Range("C6").Value = 0.9945668344998445
$ws.Range("D6").Value = 0.995538677544038
$ws.Range("E6").Value = 0.9954158728103945
$ws.Range("F6").Value = 0.9945668344998445
$ws.Range("G6").Value = 0.9948902203841825
$ws.Range("H6").Value = 0.9953592397922771
$ws.Range("I6").Value = 0.9949000558573908
$ws.Range("J6").Value = 0.995538677544038
$ws.Range("K6").Value = 0.9945668344998445
$ws.Range("L6").Value = 0.9954158728103945
$ws.Range("M6").Value = 0.9954772751772163
$ws.Range("N6").Value = 0.9954772751772163
$ws.Range("O6").Value = 0.9952815902462051
$ws.Range("P6").Value = 0.9951737949514258
$ws.Range("Q6").Value = 0.9951737949514258
$ws.Range("R6").Value = 0.9950220548385305
$ws.Range("S6").Value = 0.9950220548385305
$ws.Range("T6").Value = 0.995111816814688
$ws.Range("B7").Value = 'OffsetF'
$ws.Range("C7").Value = 1.254046386786579
$ws.Range("D7").Value = 0.5902852846375484
$ws.Range("E7").Value = 1.013491075267944
$ws.Range("F7").Value = 1.254046386786579
$ws.Range("G7").Value = 0.817393697849082
$ws.Range("H7").Value = 1.007450620122221
$ws.Range("I7").Value = 1.080894387589682
$ws.Range("J7").Value = 0.5902852846375484
$ws.Range("K7").Value = 1.254046386786579
$ws.Range("L7").Value = 1.013491075267944
$ws.Range("M7").Value = 0.8018881799527461
$ws.Range("N7").Value = 0.8018881799527461
$ws.Range("O7").Value = 0.8070566859181914
$ws.Range("P7").Value = 0.9526075822306902
$ws.Range("Q7").Value = 0.9526075822306902
$ws.Range("R7").Value = 1.027967283369662
$ws.Range("S7").Value = 1.027967283369662
$ws.Range("T7").Value = 0.9605935753755093
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 'OffsetA'
$ws.Range("C8").Value = 0.9282756955364243
$ws.Range("D8").Value = 0.9520608028659341
$ws.Range("E8").Value = 1.044576397075839
$ws.Range("F8").Value = 0.9282756955364243
$ws.Range("G8").Value = 0.9556903310574065
$ws.Range("H8").Value = 1.127367987636346
$ws.Range("I8").Value = 1.012396795706621
$ws.Range("J8").Value = 0.9520608028659341
$ws.Range("K8").Value = 0.9282756955364243
$ws.Range("L8").Value = 1.044576397075839
$ws.Range("M8").Value = 0.9983185999708863
$ws.Range("N8").Value = 0.9983185999708863
$ws.Range("O8").Value = 0.9841091769997264
$ws.Range("P8").Value = 0.974970965159399
$ws.Range("Q8").Value = 0.974970965159399
$ws.Range("R8").Value = 0.9632971477536553
$ws.Range("S8").Value = 0.9632971477536553
$ws.Range("T8").Value = 1.003394668313095
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 'RD Single'
$ws.Range("C9").Value = 1.98
$ws.Range("D9").Value = 0.22
$ws.Range("E9").Value = 0.83
$ws.Range("F9").Value = 1.98
$ws.Range("G9").Value = 0.64
$ws.Range("H9").Value = 0.6899999999999999
$ws.Range("I9").Value = 1.14
$ws.Range("J9").Value = 0.22
$ws.Range("K9").Value = 1.98
$ws.Range("L9").Value = 0.83
$ws.Range("M9").Value = 0.525
$ws.Range("N9").Value = 0.525
$ws.Range("O9").Value = 0.5633333333333334
$ws.Range("P9").Value = 1.01
$ws.Range("Q9").Value = 1.01
$ws.Range("R9").Value = 1.2525
$ws.Range("S9").Value = 1.2525
$ws.Range("T9").Value = 0.9166666666666666
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 'TD Single'
$ws.Range("C10").Value = 1.07
$ws.Range("D10").Value = 0.2
$ws.Range("E10").Value = 1.27
$ws.Range("F10").Value = 1.07
$ws.Range("G10").Value = 0.42
$ws.Range("H10").Value = 1.94
$ws.Range("I10").Value = 1.23
$ws.Range("J10").Value = 0.2
$ws.Range("K10").Value = 1.07
$ws.Range("L10").Value = 1.27
$ws.Range("M10").Value = 0.735
$ws.Range("N10").Value = 0.735
$ws.Range("O10").Value = 0.63
$ws.Range("P10").Value = 0.8466666666666667
$ws.Range("Q10").Value = 0.8466666666666667
$ws.Range("R10").Value = 0.9025000000000001
$ws.Range("S10").Value = 0.9025000000000001
$ws.Range("T10").Value = 1.021666666666667
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 'HexGrid-90degTilt5degRes'
$ws.Range("C11").Value = 0.9971325291172008
$ws.Range("D11").Value = 0.9943018859617164
$ws.Range("E11").Value = 0.9945876637685392
$ws.Range("F11").Value = 0.9971325291172008
$ws.Range("G11").Value = 0.9948649635615939
$ws.Range("H11").Value = 0.9928729087590665
$ws.Range("I11").Value = 0.9949228420563254
$ws.Range("J11").Value = 0.9943018859617164
$ws.Range("K11").Value = 0.9971325291172008
$ws.Range("L11").Value = 0.9945876637685392
$ws.Range("M11").Value = 0.9944447748651278
$ws.Range("N11").Value = 0.9944447748651278
$ws.Range("O11").Value = 0.9945848377639499
$ws.Range("P11").Value = 0.9953406929491523
$ws.Range("Q11").Value = 0.9953406929491523
$ws.Range("R11").Value = 0.9957886519911644
$ws.Range("S11").Value = 0.9957886519911644
$ws.Range("T11").Value = 0.994780465537407

# --- Apply existing header/index formatting (style index reused, no new style created) ---
$ws.Range("B1").Copy()
$ws.Range("T1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("A11").PasteSpecial(-4122)

# Clear clipboard/marching-ants selection state left over from PasteSpecial
$excel.CutCopyMode = $false
